# Update Item Name (D) / UOM (E) cell text for the rows whose underlying
# shared-string order changed, re-pairing each row with the correct
# Item Name / UOM text (values of A, B, C, F, G stay untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;  D = "Dinafex 180mg Tablet";            E = "30's" },
    @{ Row = 4;  D = "Dinafex 120mg Tablet";            E = "30's" },
    @{ Row = 7;  D = "Etorix 60mg Tablet - 40's";       E = "40's" },
    @{ Row = 8;  D = "Etorix 90mg Tablet";              E = "30's" },
    @{ Row = 9;  D = "Etorix 120mg Tablet";             E = "20's" },
    @{ Row = 11; D = "Flucloxin 500mg Capsule - 36's";  E = "36 's" },
    @{ Row = 12; D = "Flucloxin 500mg Capsule";         E = "30 's" },
    @{ Row = 15; D = "Ketonic 30mg Injection";          E = "5 's" },
    @{ Row = 16; D = "Ketonic 10mg Tablet";             E = "20's" },
    @{ Row = 18; D = "Kynol TR 100mg Capsule";          E = "50 's" },
    @{ Row = 19; D = "Kynol TR 200mg Capsule";          E = "30 's" },
    @{ Row = 25; D = "Zithrox 250mg Tablet - 6's";      E = "6's" },
    @{ Row = 26; D = "Zithrox 500mg Tablet";            E = "6 's" },
    @{ Row = 27; D = "Zithrox 30ml Dry Suspension";     E = "30ml" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
